$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B to fit the new, longer note text
$ws.Columns.Item(2).ColumnWidth = 65.75

# Add the new expense row (19/5/2018 - LoRa module parts), matching the
# formatting of the row above it
$ws.Range("A12:D12").Copy() | Out-Null
$ws.Range("A13:D13").PasteSpecial(-4122) | Out-Null
$ws.Range("A13").Value2 = 43239
$ws.Range("B13").Value2 = "Tụ 10uF + Diode 1N4007 + Nhíp kẹp linh kiện + Atmeage 328p"
$ws.Range("C13").Value2 = 223500
$ws.Range("D13").Value2 = $ws.Range("D12").Value2

$excel.CutCopyMode = 0

$ws.Range("D19").Select() | Out-Null
